$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '70.841.06'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -0.46%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.844.65'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +0.79%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '698.50'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.47%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '171.84'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -1.29%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.843.00'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.81%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("E9").Value = '  -0.74%  '
$ws.Range("E10").Value = '  -1.76%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.26'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -2.06%  '
$ws.Range("E12").Value = '  -1.16%  '
$ws.Range("E13").Value = '  -0.64%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.16'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.90%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.493.04'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +0.79%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.843.25'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +0.74%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '70.899.71'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.40%  '
$ws.Range("E18").Value = '  -1.39%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.115'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.55%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.38'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -3.55%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.73'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -4.76%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '494.26'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +1.89%  '
$ws.Range("E23").Value = '  -0.32%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.53'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.49%  '
$ws.Range("E25").Value = '  +0.72%  '
$ws.Range("B26").Value = 'InternetComputer(DFINITY)'
$ws.Range("C26").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.15'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -3.36%  '
$ws.Range("B27").Value = 'RenderToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.54'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -1.91%  '
$ws.Range("E28").Value = '  -4.66%  '
$ws.Range("E29").Value = '  -0.04%  '
$ws.Range("E30").Value = '  -1.05%  '
$ws.Range("E31").Value = '  -1.73%  '
$ws.Range("E32").Value = '  -2.70%  '
$ws.Range("B33").Value = 'Kaspa'
$ws.Range("C33").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.182'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +1.59%  '
$ws.Range("B34").Value = 'EthereumClassic'
$ws.Range("C34").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '29.41'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -1.13%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.801.14'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.95%  '
$ws.Range("E36").Value = '  -1.72%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.04%  '
$ws.Range("E38").Value = '  -0.90%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.38'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +5.58%  '
$ws.Range("E40").Value = '  +6.53%  '
$ws.Range("E41").Value = '  -0.49%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.32'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -5.48%  '
$ws.Range("E43").Value = '  +0.03%  '
$ws.Range("E44").Value = '  +0.14%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.000311'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -7.25%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '163.48'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +1.64%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '48.72'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -1.74%  '
$ws.Range("E48").Value = '  -1.21%  '
$ws.Range("E49").Value = '  +0.35%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '43.21'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -5.57%  '
$ws.Range("E51").Value = '  -5.25%  '
